# Update Excel template with additional fields for Azure storage
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing header cells (row 1) to friendlier labels -------------
$ws.Range("E1").Value = "Audit Role"
$ws.Range("F1").Value = "AgentID"
$ws.Range("G1").Value = "AgentName"
$ws.Range("H1").Value = "PBX ID"
$ws.Range("I1").Value = "Partner Name"
$ws.Range("J1").Value = "Customer Mobile"
$ws.Range("K1").Value = "Call Duration"
$ws.Range("L1").Value = "Call Type"
$ws.Range("M1").Value = "Sub Type"
$ws.Range("N1").Value = "Sub sub Type"
$ws.Range("R1").Value = "Advisor Category"

# --- Insert a new "Campaign" column before the existing callId column ------
# (shifts callId: V->W and callDate: W->X, values/types move natively so the
# forced-text "callDate" cell doesn't get re-interpreted as a serial date)
$ws.Columns("V:V").Insert()

# The insert operation drops the explicit width that used to live on the
# column it split in two (old column V / index 22); restore it so every
# pre-existing <col> definition keeps the same declared width as before.
# (Columns W/X inherit the old V/W widths automatically via the shift.)
$ws.Columns("V:V").ColumnWidth = 15

# --- Populate the new column's header and sample data ----------------------
$ws.Range("V1").Value = "Campaign"
$ws.Range("V2").Value = "Spring Promo 2025"
